# Move shape "6" (the "Draw()" rectangle) to the left, and update every
# connector that is attached to it (plus a couple of other connectors
# whose routing geometry shifted as a result), matching the target OOXML.
#
# NOTE: PowerPoint's Left/Top/Width/Height setters in this COM host are
# backed by single-precision (32-bit) floats, and EMU-from-points
# conversion truncates rather than rounds. A naive "$emu / 12700" can
# therefore land one EMU short. EmuToPt() below searches for the nearest
# point value (in the direction needed) whose float32 round-trip lands
# exactly on the requested EMU value, so every off/ext value matches the
# target precisely.

function EmuToPt([double]$emu) {
    $pt = $emu / 12700.0
    for ($i = 0; $i -lt 20000; $i++) {
        $single = [float]$pt
        $back = [math]::Floor([double]$single * 12700.0)
        if ($back -eq $emu) {
            return $pt
        }
        $pt += 0.0000001
    }
    return $pt
}

function Set-ShapeGeometry($shape, $left, $top, $width, $height) {
    if ($null -ne $left)   { $shape.Left   = EmuToPt $left }
    if ($null -ne $top)    { $shape.Top    = EmuToPt $top }
    if ($null -ne $width)  { $shape.Width  = EmuToPt $width }
    if ($null -ne $height) { $shape.Height = EmuToPt $height }
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shapesById = @{}
foreach ($shape in $s.Shapes) {
    $shapesById[$shape.Id] = $shape
}

# 1) Shape 6 ("5 Rectángulo" / Draw()) moves left.
Set-ShapeGeometry $shapesById[6] 2267744 2132856 806489 360040

# 2) Connector 42 ("41 Forma"): Grid() -> Draw(), shrinks horizontally.
Set-ShapeGeometry $shapesById[42] 1432451 1500786 1238538 632070

# 3) Connector 44 ("43 Conector angular"): Draw() -> Paint().
Set-ShapeGeometry $shapesById[44] 1929307 1391174 792088 691276
$shapesById[44].Adjustments.Item(1) = 0.5

# 4) Connector 52 ("51 Conector angular"): Paint() -> Draw().
Set-ShapeGeometry $shapesById[52] 2354154 1160748 720079 1152128
$shapesById[52].Adjustments.Item(1) = 1.31747

# 5) Connector 245 ("244 Forma"): now starts on shape 6 instead of shape 5.
$shapesById[245].ConnectorFormat.BeginConnect($shapesById[6], 3)
Set-ShapeGeometry $shapesById[245] 3074233 1624800 849695 688076
$shapesById[245].Adjustments.Item(1) = 0.5

# 6) Connector 562 ("561 Forma"): rectangle 88 -> Draw().
Set-ShapeGeometry $shapesById[562] 1920429 3243456 2088232 587112

# 7) Connector 80 ("79 Forma"): Draw() -> rectangle 77.
Set-ShapeGeometry $shapesById[80] 1403649 1801620 576064 1958617

# 8) Connector 690 ("689 Conector angular"): Setup() -> Draw().
Set-ShapeGeometry $shapesById[690] 1706081 2168860 561663 144016
